$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2628.7144
$ws.Range("I40").Value = 2531.5
$ws.Range("J40").Value = 2758.3333
$ws.Range("K40").Value = 2531.5
$ws.Range("L40").Value = 2758.3333
$ws.Range("M40").Value = -2356.5
$ws.Range("N40").Value = -3108.3333
$ws.Range("H64").Value = 3070.5
$ws.Range("I64").Value = 2997.75
$ws.Range("K64").Value = 2997.75
$ws.Range("M64").Value = -2749.75
$ws.Range("H67").Value = 3070.5
$ws.Range("I67").Value = 2997.75
$ws.Range("K67").Value = 2997.75
$ws.Range("M67").Value = -2139.75
$ws.Range("H76").Value = 3810.6667
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 4299.2
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 4299.2
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -4929.2
$ws.Range("H79").Value = 3810.6667
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 4299.2
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 4299.2
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -6483.2
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480
$ws.Range("H137").Value = 38287.37
$ws.Range("I137").Value = 884.53845
$ws.Range("K137").Value = 2653.61535
$ws.Range("M137").Value = -103.61535
$ws.Range("H138").Value = 3113.7285
$ws.Range("I138").Value = 1751.9546
$ws.Range("J138").Value = 3621.5085
$ws.Range("K138").Value = 5255.8638
$ws.Range("L138").Value = 10864.5255
$ws.Range("M138").Value = -115.8638000000001
$ws.Range("N138").Value = -21144.5255

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4134
$ws.Range("I63").Value = 3560.8
$ws.Range("K63").Value = 3560.8
$ws.Range("M63").Value = -2874.8
$ws.Range("H66").Value = 4134
$ws.Range("I66").Value = 3560.8
$ws.Range("K66").Value = 17804
$ws.Range("M66").Value = -14372
$ws.Range("H80").Value = 22600
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 24360
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 24360
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -26356
$ws.Range("H83").Value = 22600
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 24360
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 73080
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -83064

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2182.353
$ws.Range("I105").Value = 1866.6666
$ws.Range("K105").Value = 1866.6666
$ws.Range("M105").Value = -119.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 215.16667
$ws.Range("I22").Value = 296.66666
$ws.Range("K22").Value = 296.66666
$ws.Range("M22").Value = 53.33334000000002
$ws.Range("H68").Value = 24550
$ws.Range("J68").Value = 24550
$ws.Range("L68").Value = 24550
$ws.Range("N68").Value = -26048
$ws.Range("H71").Value = 24550
$ws.Range("J71").Value = 24550
$ws.Range("L71").Value = 73650
$ws.Range("N71").Value = -81138
$ws.Range("H74").Value = 19844.445
$ws.Range("J74").Value = 19844.445
$ws.Range("L74").Value = 19844.445
$ws.Range("N74").Value = -21592.445
$ws.Range("H77").Value = 19844.445
$ws.Range("J77").Value = 19844.445
$ws.Range("L77").Value = 59533.335
$ws.Range("N77").Value = -68269.33499999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 6600.5
$ws.Range("I82").Value = 606.5
$ws.Range("J82").Value = 7799.3
$ws.Range("K82").Value = 1819.5
$ws.Range("L82").Value = 23397.9
$ws.Range("M82").Value = -1413.5
$ws.Range("N82").Value = -24209.9
$ws.Range("H85").Value = 6600.5
$ws.Range("I85").Value = 606.5
$ws.Range("J85").Value = 7799.3
$ws.Range("K85").Value = 1819.5
$ws.Range("L85").Value = 23397.9
$ws.Range("M85").Value = -415.5
$ws.Range("N85").Value = -26205.9

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 5000
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("N26").Value = -5560
$ws.Range("H50").Value = 5000
$ws.Range("J50").Value = 5000
$ws.Range("L50").Value = 5000
$ws.Range("N50").Value = -5996
$ws.Range("H52").Value = 19922.389
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 19922.389
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 19922.389
$ws.Range("M52").ClearContents()  # was -2241
$ws.Range("N52").Value = -20440.389
$ws.Range("H53").Value = 7360.75
$ws.Range("J53").Value = 7360.75
$ws.Range("L53").Value = 7360.75
$ws.Range("N53").Value = -8622.75
$ws.Range("H70").Value = 75406360
$ws.Range("I70").Value = 207354740
$ws.Range("J70").Value = 7285.7144
$ws.Range("K70").Value = 207354740
$ws.Range("L70").Value = 7285.7144
$ws.Range("M70").Value = -207354470
$ws.Range("N70").Value = -7825.7144
$ws.Range("H73").Value = 75406360
$ws.Range("I73").Value = 207354740
$ws.Range("J73").Value = 7285.7144
$ws.Range("K73").Value = 207354740
$ws.Range("L73").Value = 7285.7144
$ws.Range("M73").Value = -207353804
$ws.Range("N73").Value = -9157.714400000001
$ws.Range("H80").Value = 4123.4614
$ws.Range("I80").Value = 7276.25
$ws.Range("J80").Value = 2722.2222
$ws.Range("K80").Value = 7276.25
$ws.Range("L80").Value = 2722.2222
$ws.Range("M80").Value = -6278.25
$ws.Range("N80").Value = -4718.2222
$ws.Range("H83").Value = 4123.4614
$ws.Range("I83").Value = 7276.25
$ws.Range("J83").Value = 2722.2222
$ws.Range("K83").Value = 36381.25
$ws.Range("L83").Value = 13611.111
$ws.Range("M83").Value = -31389.25
$ws.Range("N83").Value = -23595.111

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 11976.167
$ws.Range("J44").Value = 11976.167
$ws.Range("L44").Value = 11976.167
$ws.Range("N44").Value = -12888.167

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 9038
$ws.Range("J61").Value = 13057
$ws.Range("L61").Value = 13057
$ws.Range("N61").Value = -13641
